# Update COVAX_AGG package metadata for DHIS2 v2.30 (and v2.33) release.
# 1. Bump the package identifier timestamp on the "Package info" sheet.
# 2. Refresh the "Last updated" date for every indicator row on the
#    "indicators" sheet.
# 3. Correct the wording "stock on hand" -> "stock at hand" (and normalize
#    the surrounding quoting / wording) in several indicator descriptions.

$wb = $excel.ActiveWorkbook

# --- 1. Package info sheet -------------------------------------------------
$infoSheet = $wb.Worksheets.Item("Package info")
$infoSheet.Range("B7").Value = "COVAX_AGG_DASHBOARD_V1.0_DHIS2.30_2021-01-29T11:09"

# --- 2 & 3. indicators sheet ------------------------------------------------
$indSheet = $wb.Worksheets.Item("indicators")

# Rows with a corrected Description (column E) value.
$newDescription = "Opening balance equals the physical 'stock at hand count' of the previous period"

$descriptionUpdates = @{
    2   = $newDescription
    3   = $newDescription
    8   = $newDescription
    22  = "(Closing balance-Stock on hand)/Stock at hand"
    25  = $newDescription
    29  = $newDescription
    69  = $newDescription
    80  = $newDescription
    108 = $newDescription
}

foreach ($row in $descriptionUpdates.Keys) {
    $indSheet.Range("E$row").Value = $descriptionUpdates[$row]
}

# Every indicator row (2-108) gets its "Last updated" date (column I) bumped.
# The "Last updated" column stores a plain text value (e.g. "2021-01-29"),
# not a real date. A leading apostrophe forces Excel to keep it as text
# instead of auto-converting the date-like string into a date serial
# number (which would also silently swap in a date number format).
for ($row = 2; $row -le 108; $row++) {
    $indSheet.Range("I$row").Value = "'2021-01-29"
}
